$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 45826.043
$ws.Range("J17").Value = 45826.043
$ws.Range("L17").Value = 137478.129
$ws.Range("N17").Value = -137814.129

$ws.Range("H70").Value = 3209.7778
$ws.Range("I70").Value = 1996
$ws.Range("J70").Value = 3361.5
$ws.Range("K70").Value = 5988
$ws.Range("L70").Value = 10084.5
$ws.Range("M70").Value = -5718
$ws.Range("N70").Value = -10624.5

$ws.Range("H73").Value = 3209.7778
$ws.Range("I73").Value = 1996
$ws.Range("J73").Value = 3361.5
$ws.Range("K73").Value = 5988
$ws.Range("L73").Value = 10084.5
$ws.Range("M73").Value = -5052
$ws.Range("N73").Value = -11956.5

$ws.Range("H76").Value = 2593.5
$ws.Range("I76").Value = 2593.5
$ws.Range("K76").Value = 2593.5
$ws.Range("M76").Value = -2278.5

$ws.Range("H79").Value = 2593.5
$ws.Range("I79").Value = 2593.5
$ws.Range("K79").Value = 2593.5
$ws.Range("M79").Value = -1501.5

$ws.Range("H111").Value = 1774.3334
$ws.Range("I111").Value = 1272.4445
$ws.Range("J111").Value = 3280
$ws.Range("K111").Value = 3817.3335
$ws.Range("L111").Value = 9840
$ws.Range("M111").Value = -750.3335000000002
$ws.Range("N111").Value = -15974

$ws.Range("H135").Value = 4535.0713
$ws.Range("I135").Value = 4653.154
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 41878.38600000001
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -39343.38600000001
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 32279
$ws.Range("I3").Value = 7334.3335
$ws.Range("K3").Value = 7334.3335
$ws.Range("M3").Value = -7219.3335

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H16").Value = 829.3333
$ws.Range("I16").Value = 494
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 494
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -207
$ws.Range("N16").Value = -2074

$ws.Range("H32").Value = 231915.36
$ws.Range("I32").Value = 234982.84
$ws.Range("K32").Value = 234982.84
$ws.Range("M32").Value = -234695.84

$ws.Range("H36").Value = 3612.5
$ws.Range("I36").Value = 3612.5
$ws.Range("K36").Value = 3612.5
$ws.Range("M36").Value = -3266.5

$ws.Range("H74").Value = 5134.136
$ws.Range("I74").Value = 2938.8333
$ws.Range("K74").Value = 2938.8333
$ws.Range("M74").Value = -2064.8333

$ws.Range("H77").Value = 5134.136
$ws.Range("I77").Value = 2938.8333
$ws.Range("K77").Value = 14694.1665
$ws.Range("M77").Value = -10326.1665

$ws.Range("H108").Value = 57699.4
$ws.Range("J108").Value = 57699.4
$ws.Range("L108").Value = 57699.4
$ws.Range("N108").Value = -65379.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1710.1111
$ws.Range("I5").Value = 270.14285
$ws.Range("J5").Value = 6750
$ws.Range("K5").Value = 270.14285
$ws.Range("L5").Value = 6750
$ws.Range("M5").Value = -157.14285
$ws.Range("N5").Value = -6976

$ws.Range("H20").Value = 1345.2106
$ws.Range("J20").Value = 1375
$ws.Range("L20").Value = 1375
$ws.Range("N20").Value = -1869

$ws.Range("H86").Value = 5723.875
$ws.Range("I86").Value = 2542.7144
$ws.Range("K86").Value = 2542.7144
$ws.Range("M86").Value = -1419.7144

$ws.Range("H89").Value = 5723.875
$ws.Range("I89").Value = 2542.7144
$ws.Range("K89").Value = 12713.572
$ws.Range("M89").Value = -7097.572

$ws.Range("H94").Value = 1368.6511
$ws.Range("I94").Value = 1410.3529
$ws.Range("J94").Value = 1211.1111
$ws.Range("K94").Value = 1410.3529
$ws.Range("L94").Value = 1211.1111
$ws.Range("M94").Value = -959.3529000000001
$ws.Range("N94").Value = -2113.1111

$ws.Range("H134").Value = 2792.0637
$ws.Range("I134").Value = 1568.5358
$ws.Range("K134").Value = 4705.607400000001
$ws.Range("M134").Value = -2170.607400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1372.0834
$ws.Range("I22").Value = 783.5
$ws.Range("J22").Value = 1960.6666
$ws.Range("K22").Value = 783.5
$ws.Range("L22").Value = 1960.6666
$ws.Range("M22").Value = -433.5
$ws.Range("N22").Value = -2660.6666

$ws.Range("H58").Value = 2684.1667
$ws.Range("I58").Value = 2900
$ws.Range("J58").Value = 2612.2222
$ws.Range("K58").Value = 2900
$ws.Range("L58").Value = 2612.2222
$ws.Range("M58").Value = -2697
$ws.Range("N58").Value = -3018.2222

$ws.Range("H105").Value = 2722.2
$ws.Range("I105").Value = 1203.6666
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 1203.6666
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = 543.3334
$ws.Range("N105").Value = -8494

$ws.Range("H122").Value = 2554.9395
$ws.Range("I122").Value = 2721.45
$ws.Range("J122").Value = 2298.7693
$ws.Range("K122").Value = 8164.349999999999
$ws.Range("L122").Value = 6896.3079
$ws.Range("M122").Value = -5714.349999999999
$ws.Range("N122").Value = -11796.3079

$ws.Range("H132").Value = 5228.222
$ws.Range("I132").Value = 6222.2
$ws.Range("K132").Value = 18666.6
$ws.Range("M132").Value = -16136.6

$ws.Range("H134").Value = 2147.3447
$ws.Range("I134").Value = 1443
$ws.Range("K134").Value = 4329
$ws.Range("M134").Value = -1794

$ws.Range("H136").Value = 2684.1667
$ws.Range("I136").Value = 2900
$ws.Range("J136").Value = 2612.2222
$ws.Range("K136").Value = 8700
$ws.Range("L136").Value = 7836.6666
$ws.Range("M136").Value = -6150
$ws.Range("N136").Value = -12936.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 484.21738
$ws.Range("I2").Value = 212
$ws.Range("J2").Value = 994.625
$ws.Range("K2").Value = 1272
$ws.Range("L2").Value = 5967.75
$ws.Range("M2").Value = -1159
$ws.Range("N2").Value = -6193.75

$ws.Range("H38").Value = 644.8261
$ws.Range("I38").Value = 217.27272
$ws.Range("J38").Value = 1036.75
$ws.Range("K38").Value = 651.81816
$ws.Range("L38").Value = 3110.25
$ws.Range("M38").Value = -304.81816
$ws.Range("N38").Value = -3804.25

$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 60000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -60430

$ws.Range("H131").Value = 3833044.5
$ws.Range("I131").Value = 11366274
$ws.Range("J131").Value = 66429.56
$ws.Range("K131").Value = 34098822
$ws.Range("L131").Value = 199288.68
$ws.Range("M131").Value = -34093782
$ws.Range("N131").Value = -209368.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 250450
$ws.Range("I14").Value = 250450
$ws.Range("K14").Value = 250450
$ws.Range("M14").Value = -250282

$ws.Range("H62").Value = 46962.25
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 46962.25
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 46962.25
$ws.Range("N62").Value = -48334.25
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 46962.25
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 46962.25
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 140886.75
$ws.Range("N65").Value = -147750.75
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 3015.8262
$ws.Range("I132").Value = 2469.3333
$ws.Range("K132").Value = 7407.999899999999
$ws.Range("M132").Value = -4877.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1968.6
$ws.Range("J12").Value = 1968.6
$ws.Range("L12").Value = 1968.6
$ws.Range("N12").Value = -2308.6

$ws.Range("H46").Value = 6109.4736
$ws.Range("I46").Value = 20618.6
$ws.Range("J46").Value = 927.6429000000001
$ws.Range("K46").Value = 20618.6
$ws.Range("L46").Value = 927.6429000000001
$ws.Range("M46").Value = -20430.6
$ws.Range("N46").Value = -1303.6429

$ws.Range("H61").Value = 12502679
$ws.Range("I61").Value = 20002292
$ws.Range("J61").Value = 3323.5
$ws.Range("K61").Value = 20002292
$ws.Range("L61").Value = 3323.5
$ws.Range("M61").Value = -20002090
$ws.Range("N61").Value = -3727.5

$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248

$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240

$ws.Range("H113").Value = 12502679
$ws.Range("I113").Value = 20002292
$ws.Range("J113").Value = 3323.5
$ws.Range("K113").Value = 20002292
$ws.Range("L113").Value = 3323.5
$ws.Range("M113").Value = -20000122
$ws.Range("N113").Value = -7663.5

$ws.Range("H136").Value = 4463.615
$ws.Range("I136").Value = 4239.125
$ws.Range("K136").Value = 12717.375
$ws.Range("M136").Value = -10167.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 503.4
$ws.Range("I113").Value = 415.8889
$ws.Range("K113").Value = 1247.6667
$ws.Range("M113").Value = 922.3333

$ws.Range("H122").Value = 10003.286
$ws.Range("I122").Value = 10238.823
$ws.Range("K122").Value = 30716.469
$ws.Range("M122").Value = -28266.469

$ws.Range("H136").Value = 25911.293
$ws.Range("I136").Value = 38027.555
$ws.Range("K136").Value = 114082.665
$ws.Range("M136").Value = -111532.665
